# Fix typos in the "Geoscience Australia profile" example diagram (slide 4).
# The four "Base Specification"/"format" rounded rectangles had their labels
# mixed up; restore the correct XSD / XML Schema Definition / Schematron / XML
# wording.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$s.Shapes.Item("Rounded Rectangle 53").TextFrame.TextRange.Text = "XML"
$s.Shapes.Item("Rounded Rectangle 61").TextFrame.TextRange.Text = "XML Schema Definition"
$s.Shapes.Item("Rounded Rectangle 41").TextFrame.TextRange.Text = "Schematron"
$s.Shapes.Item("Rounded Rectangle 56").TextFrame.TextRange.Text = "XML"
